$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185:198 down to 186:199.
$ws.Range("A185").EntireRow.Insert()

# Populate the newly inserted row 185 with the new price-observation record.
$ws.Range("A185").Value = 1
$ws.Range("B185").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C185").Value = "Arica y Parinacota"
$ws.Range("D185").Value = 44610
$ws.Range("E185").Value = 15
$ws.Range("F185").Value = "Fruta"
$ws.Range("G185").Value = 100102
$ws.Range("H185").Value = "Cítricos"
$ws.Range("I185").Value = 100102003
$ws.Range("J185").Value = "Limón"
$ws.Range("K185").Value = "Tahití"
$ws.Range("L185").Value = "Primera"
$ws.Range("M185").Value = 300
$ws.Range("N185").Value = 35000
$ws.Range("O185").Value = 36000
$ws.Range("P185").Value = 35500
$ws.Range("Q185").Value = "`$/caja 24 kilos"
$ws.Range("R185").Value = "Perú"
$ws.Range("S185").Value = 1479
$ws.Range("T185").Value = 24
